$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: the =TODAY() formula is replaced with the static date value it had
# last resolved to, so the cell no longer recalculates.
$ws.Range("A3").Value = 45866

# Duplicate row 3's formatting onto the two new rows (4 and 5) before
# filling in their content, same as copy/pasting the row and editing it.
$ws.Range("A3:E3").Copy()
$ws.Range("A4:E4").PasteSpecial(-4122)
$ws.Range("A5:E5").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 4 - Tuesday
$ws.Range("A4").Value = 45867
$ws.Range("B4").Value = "Tuesday"
$ws.Range("C4").Value = " •full stack web development "
$ws.Range("D4").Value = "• Practiced Node.js concepts and setup.`n• Understood and implemented MongoDB connection using Mongoose.`n• Explored MongoDB Atlas and created a new cluster for remote database access.`n• Understood each line of the MongoDB connection logic and its error handling mechanism."

# Row 5 - Wednesday
$ws.Range("A5").Value = 45868
$ws.Range("B5").Value = "Wednesday"
$ws.Range("C5").Value = " •full stack web development "
$ws.Range("D5").Value = "• Practiced on creating servers with express js`n•worked on  making a basic photo gallery app by using node js`n•Practiced on writing Javascript syntaxes"

# Row heights (reflecting the wrapped-text reflow for the new content)
$ws.Range("A3:E3").RowHeight = 202.9
$ws.Range("A4:E4").RowHeight = 162
$ws.Range("A5:E5").RowHeight = 81

# Column widths (auto-fit side effect of the new content)
$ws.Columns.Item(1).ColumnWidth = 9.8337
$ws.Columns.Item(2).ColumnWidth = 11.0003
$ws.Columns.Item(3).ColumnWidth = 42.8337
$ws.Columns.Item(4).ColumnWidth = 58.3337
$ws.Columns.Item(5).ColumnWidth = 41.8337

# Update the view: scroll so row 3 is near the top, with E4:E5 selected
$excel.ActiveWindow.ScrollRow = 3
$ws.Range("E4:E5").Select()
